# Fruta / hortaliza, semanal
# Insert a new weekly price-observation row above row 108 (pushing the
# existing rows 108-227 down to 109-228) and populate the new row with the
# latest week's data for "Acelga" @ Terminal Hortofrutícola Agro Chillán.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 108..227 down to 109..228.
$ws.Rows("108").Insert()

# Fill the newly inserted row 108 with the new weekly observation.
$ws.Cells.Item(108, 1).Value = 7
$ws.Cells.Item(108, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(108, 3).Value = "Ñuble"
$ws.Cells.Item(108, 4).Value = 44664
$ws.Cells.Item(108, 5).Value = 16
$ws.Cells.Item(108, 6).Value = 100112009
$ws.Cells.Item(108, 7).Value = "Acelga"
$ws.Cells.Item(108, 8).Value = "Sin especificar"
$ws.Cells.Item(108, 9).Value = "Primera"
$ws.Cells.Item(108, 10).Value = 200
$ws.Cells.Item(108, 11).Value = 550
$ws.Cells.Item(108, 12).Value = 600
$ws.Cells.Item(108, 13).Value = 575
$ws.Cells.Item(108, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(108, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(108, 16).Value = 575
$ws.Cells.Item(108, 17).Value = 1
$ws.Cells.Item(108, 18).Value = "Hortaliza"
